$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7591879963874817
$ws.Range("B1").Value = 0.9178144335746765
$ws.Range("C1").Value = 1.24742865562439
$ws.Range("D1").Value = 2.941295862197876
$ws.Range("E1").Value = 2.589521408081055
